# Insert a new data row above row 500 (shifts existing rows 500:534 down to 501:535)
# and populate the new row 500 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 500 and below down by one row, inserting a blank row at 500.
$ws.Rows.Item(500).Insert()

# Populate the newly inserted row 500.
$ws.Cells.Item(500, 1).Value = 4
$ws.Cells.Item(500, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(500, 3).Value = "Los Lagos"
$ws.Cells.Item(500, 4).Value = 45013
$ws.Cells.Item(500, 5).Value = 10
$ws.Cells.Item(500, 6).Value = 100114013
$ws.Cells.Item(500, 7).Value = "Zanahoria"
$ws.Cells.Item(500, 8).Value = "Sin especificar"
$ws.Cells.Item(500, 9).Value = "Primera"
$ws.Cells.Item(500, 10).Value = 750
$ws.Cells.Item(500, 11).Value = 9000
$ws.Cells.Item(500, 12).Value = 9000
$ws.Cells.Item(500, 13).Value = 9000
$ws.Cells.Item(500, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(500, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(500, 16).Value = 450
$ws.Cells.Item(500, 17).Value = 20
$ws.Cells.Item(500, 18).Value = "Hortaliza"
